# Update the "Assigment Assessment" sheet (4th assessment column, col G)
# with the grading entries, matching the commit's recorded changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Assigment Assessment")
$ws.Activate() | Out-Null

# --- Column G (4th assignment score) for rows 3-49 ---
$gValues = @(
    90, 80, 60, 85, 90, 80, 80, 80, 90, 90, 85, 60,
    80, 90, 90, 90, 90, 80, 90, 80, 75, 90, 70,
    75, 80, 90, 60, 90, 80, 75, 60, 90, 80, 70,
    80, 85, 70, 80, 80, 80, 85, 85, 90, 80, 80,
    80, 85
)
$row = 3
foreach ($v in $gValues) {
    $ws.Cells.Item($row, 7).Value = $v
    $row = $row + 1
}

# --- Column F (3rd assignment score) for rows 42-48, previously blank ---
$fValues = @{
    42 = 90
    43 = 85
    44 = 85
    45 = 90
    46 = 80
    47 = 90
    48 = 60
}
foreach ($r in $fValues.Keys) {
    $ws.Cells.Item($r, 6).Value = $fValues[$r]
}

# --- G5 is additionally flagged: red, centered text (new distinct style) ---
$g5 = $ws.Range("G5")
$g5.Value = 60
$g5.Font.Color = 255
$g5.Font.Name = "等线"
$g5.Font.Size = 11
$g5.HorizontalAlignment = -4108
$g5.VerticalAlignment = -4108

# --- Update the last active selection, as recorded after the edit ---
$ws.Range("G48").Select() | Out-Null
